$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Friction values")

# Update µStatic values for gravel (row 7) and sand (row 8) to match µDynamic
$ws.Range("C7").Value = 0.35
$ws.Range("C8").Value = 0.3

# Update the active selection to reflect the last-edited cell
$ws.Range("F7").Select()
